$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 10 ---- (order matches original authoring sequence so that the
# shared-string table is rebuilt with the same index order as the source)
$ws.Range("D10").Value = "Trinity Doctorate Award"
$ws.Range("B10").Value = "122K, PI"
$ws.Range("A10").Value = "Trinity Doctorate Award 2024-2025, PI-led application"
$ws.Range("C10").Value = "2024-2028"
$ws.Range("E10").Value = "Mapping Macrophage Cell State Transition in Inflammation and Infection at Single Cell Resolution"

$ws.Range("E10").Font.Name = "Calibri"
$ws.Range("E10").Font.Size = 10.5
$ws.Range("E10").Font.Color = 2368548
$ws.Range("E10").VerticalAlignment = -4108

# ---- Row 11 ----
$ws.Rows.Item(11).RowHeight = 20

$ws.Range("B11").Value = "30.2K, PI"
$ws.Range("A11").Value = "Higher Education Research Equipment Grant"
$ws.Range("D11").Value = "Higher Education Authority "
$ws.Range("C11").Value = "2024-2025"
$ws.Range("E11").Value = "Cell counting solution to support TCD Omics"

$ws.Range("A11").Font.Name = "Arial"
$ws.Range("A11").Font.Size = 14
$ws.Range("A11").Font.Color = 4802889

$ws.Range("E11").Font.Name = "Calibri"
$ws.Range("E11").Font.Size = 15
$ws.Range("E11").Font.Color = 2171169

# ---- selection, matches end-state cursor ----
$ws.Range("A11").Select()
